# Fruta / hortaliza, semanal
#
# Two new daily price records (rows) are inserted into the "Femacal de La
# Calera - Limón" data range, right before the existing row that used to
# sit at sheet row 540. All subsequent rows (old 540:633) shift down by two
# rows (to 542:635), which also grows the used range from A1:T633 to
# A1:T635.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 540, pushing the old rows 540:633 down to 542:635.
$ws.Range("A540:A541").EntireRow.Insert()

# --- New row 540: Limón "1a amarillo" record dated 2021-11-04 (44504) ---
$ws.Cells.Item(540, 1).Value  = 3
$ws.Cells.Item(540, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(540, 3).Value  = "Coquimbo"
$ws.Cells.Item(540, 4).Value  = 44504
$ws.Cells.Item(540, 5).Value  = 5
$ws.Cells.Item(540, 6).Value  = "Fruta"
$ws.Cells.Item(540, 7).Value  = 100102
$ws.Cells.Item(540, 8).Value  = "Cítricos"
$ws.Cells.Item(540, 9).Value  = 100102003
$ws.Cells.Item(540, 10).Value = "Limón"
$ws.Cells.Item(540, 11).Value = "Sin especificar"
$ws.Cells.Item(540, 12).Value = "1a amarillo"
$ws.Cells.Item(540, 13).Value = 328
$ws.Cells.Item(540, 14).Value = 4500
$ws.Cells.Item(540, 15).Value = 5500
$ws.Cells.Item(540, 16).Value = 5018
$ws.Cells.Item(540, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(540, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(540, 19).Value = 314
$ws.Cells.Item(540, 20).Value = 16

# --- New row 541: Limón "2a amarillo" record dated 2021-11-04 (44504) ---
$ws.Cells.Item(541, 1).Value  = 3
$ws.Cells.Item(541, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(541, 3).Value  = "Coquimbo"
$ws.Cells.Item(541, 4).Value  = 44504
$ws.Cells.Item(541, 5).Value  = 5
$ws.Cells.Item(541, 6).Value  = "Fruta"
$ws.Cells.Item(541, 7).Value  = 100102
$ws.Cells.Item(541, 8).Value  = "Cítricos"
$ws.Cells.Item(541, 9).Value  = 100102003
$ws.Cells.Item(541, 10).Value = "Limón"
$ws.Cells.Item(541, 11).Value = "Sin especificar"
$ws.Cells.Item(541, 12).Value = "2a amarillo"
$ws.Cells.Item(541, 13).Value = 299
$ws.Cells.Item(541, 14).Value = 3500
$ws.Cells.Item(541, 15).Value = 4000
$ws.Cells.Item(541, 16).Value = 3780
$ws.Cells.Item(541, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(541, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(541, 19).Value = 236
$ws.Cells.Item(541, 20).Value = 16
